# This script turns the "valid" complex-charting test workbook into the
# "invalid" variant used by programs-charting-complex-multiple-invalid.xlsx.
#
# Summary of the required edits:
#  - Metadata sheet: rows 8/9 (Test Chart One / Test Chart Two survey rows)
#    lose their special formatting (style index 2 -> 1).
#  - Core sheet: the "code" column (A2:A5) is overwritten so it duplicates
#    the "name" column (B2:B5) instead of using the real testchartcorecode*
#    codes; the stray alternate formatting on R5 is cleared (style 2 -> 1).
#  - Test Chart One / Test Chart Two sheets: the first question's code
#    ("testchartcode0") is renamed to "PatientChartingDate" and picks up
#    the alternate formatting (style 1 -> 2).
#  - Test Chart Two sheet: the alternate formatting previously on the
#    "Question A"/"Question B" cells (C3:D3, C4:D4) is cleared (style 2 -> 1).
#
# Once the now-unused "testchartcorecode*" shared strings are no longer
# referenced by any cell, the engine recomputes sharedStrings.xml on save
# and drops them automatically (uniqueCount 69 -> 65), so we only need to
# change the cell values/styles themselves.

$wb = $excel.ActiveWorkbook
$wsMetadata = $wb.Worksheets.Item(1)   # Metadata
$wsCore     = $wb.Worksheets.Item(2)   # Core
$wsChartOne = $wb.Worksheets.Item(3)   # Test Chart One
$wsChartTwo = $wb.Worksheets.Item(4)   # Test Chart Two

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Capture the two distinct cell formats that are in play, from cells
#     that still carry them untouched at this point in the script. ---

# "Alternate" style (s=2 in the original file) - currently on Test Chart Two!C3
$wsChartTwo.Range("C3").Copy()
$wsChartOne.Range("A2").PasteSpecial($xlPasteFormats)
$wsChartTwo.Range("A2").PasteSpecial($xlPasteFormats)

# Now give the renamed "code" cells their new value.
$wsChartOne.Range("A2").Value = "PatientChartingDate"
$wsChartTwo.Range("A2").Value = "PatientChartingDate"

# "Normal" style (s=1 in the original file) - a cell that is never touched
# by this edit, e.g. the Core sheet header.
$wsCore.Range("A1").Copy()
$wsMetadata.Range("A8").PasteSpecial($xlPasteFormats)
$wsMetadata.Range("B8").PasteSpecial($xlPasteFormats)
$wsMetadata.Range("A9").PasteSpecial($xlPasteFormats)
$wsMetadata.Range("B9").PasteSpecial($xlPasteFormats)
$wsCore.Range("R5").PasteSpecial($xlPasteFormats)
$wsChartTwo.Range("C3").PasteSpecial($xlPasteFormats)
$wsChartTwo.Range("D3").PasteSpecial($xlPasteFormats)
$wsChartTwo.Range("C4").PasteSpecial($xlPasteFormats)
$wsChartTwo.Range("D4").PasteSpecial($xlPasteFormats)

# --- Update the Core sheet's "code" column so it matches "name" column ---
$wsCore.Range("A2").Value = "ComplexChartInstanceName"
$wsCore.Range("A3").Value = "ComplexChartDate"
$wsCore.Range("A4").Value = "ComplexChartType"
$wsCore.Range("A5").Value = "ComplexChartSubtype"
